$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.494.26'
$ws.Range('E2').Value = '  +0.08%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.825.94'
$ws.Range('E3').Value = '  -0.06%  '
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '316.36'
$ws.Range('E5').Value = '  +0.36%  '
$ws.Range('E6').Value = '  +0.24%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5182'
$ws.Range('E7').Value = '  +2.32%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3870'
$ws.Range('E8').Value = '  -0.97%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08297'
$ws.Range('E9').Value = '  +8.24%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.124'
$ws.Range('E10').Value = '  +1.41%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '41.92'
$ws.Range('E11').Value = '  +0.11%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.382'
$ws.Range('E12').Value = '  +1.53%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.18'
$ws.Range('E13').Value = '  +0.42%  '
$ws.Range('E14').Value = '  +0.24%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.492'
$ws.Range('E15').Value = '  -1.22%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.830.60'
$ws.Range('E16').Value = '  +0.40%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '94.01'
$ws.Range('E17').Value = '  +1.09%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001123'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06635'
$ws.Range('E19').Value = '  -0.58%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.84'
$ws.Range('E20').Value = '  +0.93%  '
$ws.Range('E21').Value = '  +0.23%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.068'
$ws.Range('E22').Value = '  -1.18%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.528.76'
$ws.Range('E23').Value = '  +0.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.49'
$ws.Range('E24').Value = '  +3.23%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.244'
$ws.Range('E25').Value = '  -0.43%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '21.12'
$ws.Range('E26').Value = '  +2.52%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '159.74'
$ws.Range('E27').Value = '  +1.97%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.037.93'
$ws.Range('E28').Value = '  +0.19%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.415'
$ws.Range('E29').Value = '  +0.90%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '125.99'
$ws.Range('E30').Value = '  +0.64%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.1099'
$ws.Range('E31').Value = '  +1.50%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.097'
$ws.Range('E32').Value = '  -2.60%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.740'
$ws.Range('E33').Value = '  +1.08%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.07595'
$ws.Range('E34').Value = '  +7.95%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.682'
$ws.Range('E35').Value = '  +0.57%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.2232'
$ws.Range('E36').Value = '  +0.21%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02372'
$ws.Range('E37').Value = '  +2.27%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.261'
$ws.Range('E38').Value = '  +2.47%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '12.06'
$ws.Range('E39').Value = '  +7.39%  '
$ws.Range('E40').Value = '  -2.22%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6415'
$ws.Range('E41').Value = '  +2.63%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.191'
$ws.Range('E42').Value = '  +0.88%  '
$ws.Range('E43').Value = '  +0.15%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.65'
$ws.Range('E44').Value = '  +1.66%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6216'
$ws.Range('E45').Value = '  +5.39%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.799'
$ws.Range('E46').Value = '  +2.13%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '128.03'
$ws.Range('E47').Value = '  +2.76%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.006'
$ws.Range('E48').Value = '  +1.30%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.205'
$ws.Range('E49').Value = '  +1.12%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06972'
$ws.Range('E50').Value = '  +0.80%  '
$ws.Range('E51').Value = '  +1.04%  '
